$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '52.312.82'
$ws.Range('E2').Value = '  +1.02%  '
$ws.Range('D3').Value = '3.014.93'
$ws.Range('E3').Value = '  +2.45%  '
$ws.Range('E4').Value = '  -0.03%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '355.02'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +0.91%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '108.43'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -2.70%  '
$ws.Range('E7').Value = '  +0.48%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -1.68%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '38.53'
$c.Style = "Normal"
$ws.Range('E10').Value = '  -2.39%  '
$ws.Range('E11').Value = '  +2.12%  '
$ws.Range('E12').Value = '  -3.80%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '19.29'
$c.Style = "Normal"
$ws.Range('E13').Value = '  -2.69%  '
$ws.Range('D14').Value = '3.489.05'
$ws.Range('E14').Value = '  +2.32%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '7.73'
$c.Style = "Normal"
$ws.Range('E15').Value = '  -3.72%  '
$ws.Range('D16').Value = '3.005.27'
$ws.Range('E16').Value = '  +1.79%  '
$ws.Range('E17').Value = '  +2.98%  '
$ws.Range('D18').Value = '52.345.54'
$ws.Range('E18').Value = '  +0.86%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '3.55'
$c.Style = "Normal"
$ws.Range('E19').Value = '  +8.26%  '
$ws.Range('E20').Value = '  -1.85%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '13.75'
$c.Style = "Normal"
$ws.Range('E21').Value = '  -5.47%  '
$ws.Range('D22').Value = '0.0₃0976'
$ws.Range('E22').Value = '  -1.37%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '69.66'
$c.Style = "Normal"
$ws.Range('E23').Value = '  -2.39%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '265.39'
$c.Style = "Normal"
$ws.Range('E24').Value = '  -2.82%  '
$ws.Range('E25').Value = '  -1.38%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '0.180'
$c.Style = "Normal"
$ws.Range('E26').Value = '  -2.07%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '7.74'
$c.Style = "Normal"
$ws.Range('E27').Value = '  +4.26%  '
$ws.Range('E28').Value = '  -1.30%  '
$ws.Range('E29').Value = '  -0.06%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '0.107'
$c.Style = "Normal"
$ws.Range('E30').Value = '  -3.84%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '6.47'
$c.Style = "Normal"
$ws.Range('E31').Value = '  +2.05%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '10.36'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -3.59%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '36.50'
$c.Style = "Normal"
$ws.Range('E33').Value = '  -2.83%  '
$ws.Range('E34').Value = '  +20.90%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '50.96'
$c.Style = "Normal"
$ws.Range('E35').Value = '  -3.64%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '0.0444'
$c.Style = "Normal"
$ws.Range('E36').Value = '  -0.93%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('E38').Value = '  -4.57%  '
$ws.Range('E40').Value = '  -4.16%  '
$ws.Range('E41').Value = '  +2.60%  '
$ws.Range('E42').Value = '  -0.70%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '22.98'
$c.Style = "Normal"
$ws.Range('E43').Value = '  -3.08%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '123.64'
$c.Style = "Normal"
$ws.Range('E44').Value = '  +9.30%  '
$ws.Range('E45').Value = '  -0.87%  '
$ws.Range('D46').Value = '2.130.27'
$ws.Range('E46').Value = '  -1.33%  '
$ws.Range('E47').Value = '  -3.74%  '
$ws.Range('E48').Value = '  -5.65%  '
$ws.Range('D49').Value = '3.310.22'
$ws.Range('E49').Value = '  +2.20%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '0.250'
$c.Style = "Normal"
$ws.Range('E50').Value = '  +1.75%  '
